# Quarterly indexing esoteric bug-fix operation
# Column A holds a date per row that marked the 1st of a quarter-start month.
# The fix re-indexes each of those dates to the 15th of the following month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2
    if ($oldVal -ne $null) {
        $serial = [double]$oldVal
        $epochDays = [math]::Floor($serial)

        # Decompose the Excel 1900-date-system serial number into y/m/d.
        $daysInMonth = @(31,28,31,30,31,30,31,31,30,31,30,31)
        $year = 1900
        $remaining = $epochDays - 1   # serial 1 == 1900-01-01
        while ($true) {
            $isLeap = (($year % 4 -eq 0 -and $year % 100 -ne 0) -or ($year % 400 -eq 0)) -or ($year -eq 1900)
            $yearDays = 365
            if ($isLeap) { $yearDays = 366 }
            if ($remaining -lt $yearDays) { break }
            $remaining = $remaining - $yearDays
            $year = $year + 1
        }
        $isLeap = (($year % 4 -eq 0 -and $year % 100 -ne 0) -or ($year % 400 -eq 0)) -or ($year -eq 1900)
        $month = 1
        for ($mi = 0; $mi -lt 12; $mi++) {
            $dim = $daysInMonth[$mi]
            if ($mi -eq 1 -and $isLeap) { $dim = 29 }
            if ($remaining -lt $dim) { $month = $mi + 1; break }
            $remaining = $remaining - $dim
        }

        $targetMonth = $month + 1
        $targetYear = $year
        if ($targetMonth -gt 12) {
            $targetMonth = $targetMonth - 12
            $targetYear = $targetYear + 1
        }

        # Recompose serial number for (targetYear, targetMonth, 15).
        $newSerial = 1
        for ($y = 1900; $y -lt $targetYear; $y++) {
            $yLeap = (($y % 4 -eq 0 -and $y % 100 -ne 0) -or ($y % 400 -eq 0)) -or ($y -eq 1900)
            if ($yLeap) { $newSerial = $newSerial + 366 } else { $newSerial = $newSerial + 365 }
        }
        $tLeap = (($targetYear % 4 -eq 0 -and $targetYear % 100 -ne 0) -or ($targetYear % 400 -eq 0)) -or ($targetYear -eq 1900)
        for ($mi = 1; $mi -lt $targetMonth; $mi++) {
            $dim = $daysInMonth[$mi - 1]
            if ($mi -eq 2 -and $tLeap) { $dim = 29 }
            $newSerial = $newSerial + $dim
        }
        $newSerial = $newSerial + (15 - 1)

        $cell.Value2 = $newSerial
    }
}
